$d = $word.ActiveDocument

# --- 1. Remove the "_GoBack" bookmark left over from the previous edit session ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Split "Consistent, effective and creative use of spelling; no errors in
#        mechanics." into "Consistent and effective use of spelling; no errors
#        in mechanics." (drop ", ... creative") while keeping the sentence split
#        across three runs, as in the authored edit. ---
$r = $d.Content
$found = $r.Find.Execute("Consistent, effective and creative use of spelling; no errors in mechanics.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $s = $r.Start

    # Replace the whole sentence with the corrected wording first.
    $r.Text = "Consistent and effective use of spelling; no errors in mechanics."

    $len1 = 10  # "Consistent"
    $len2 = 5   # " and "
    $len3 = 50  # "effective use of spelling; no errors in mechanics."

    $r1 = $d.Range($s, $s + $len1)
    $r2 = $d.Range($s + $len1, $s + $len1 + $len2)
    $r3 = $d.Range($s + $len1 + $len2, $s + $len1 + $len2 + $len3)

    # Nudge formatting on the middle run (and revert) so Word keeps the three
    # pieces as distinct runs instead of silently re-merging them.
    $r2.Font.Bold = 1
    $r2.Font.Bold = 0
}
